{"js": "\n// 1. Title\n{\n  const body = context.document.body;\n  const res = body.search(\"Galaxies: Cosmic Tapestry of Wonders\", {matchCase: true});\n  res.load(\"items\");\n  await context.sync();\n  res.items[0].insertText(\"Exploring Chemistry - The Building Blocks of Life\", \"Replace\");\n  await context.sync();\n}\n\n// 2. Author\n{\n  const body = context.document.body;\n  const res = body.search(\"Katie Leonard\", {matchCase: true});\n  res.load(\"items\");\n  await context.sync();\n  res.items[0].insertText(\"Maya Singh\", \"Replace\");\n  await context.sync();\n}\n\n// 3. Email local-part\n{\n  const body = context.document.body;\n  const res = body.search(\"intelligentwriter349@digitalworld\", {matchCase: true});\n  res.load(\"items\");\n  await context.sync();\n  res.items[0].insertText(\"mrssingh@sunnyvalehigh\", \"Replace\");\n  await context.sync();\n}\n\n// 4. Email domain suffix\n{\n  const body = context.document.body;\n  const res = body.search(\"tech\", {matchCase: true});\n  res.load(\"items\");\n  await context.sync();\n  res.items[0].insertText(\"edu\", \"Replace\");\n  await context.sync();\n}\n\n// 5. Main body paragraph (index 4) - full replacement, preserves run formatting\n{\n  const paras = context.document.body.paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n  const p = paras.items[4];\n  const full = p.getRange();\n  full.insertText(\"Chemistry is often regarded as the central science due to its profound influence on various fields, spanning biology, medicine, and materials science. Its study unveils the fundamental principles governing the behavior of matter at the atomic and molecular levels. By delving into the intricate dance of chemical reactions, we uncover the secrets behind the formation and transformation of substances, unlocking the enigmas that shape our world. Chemistry reveals the tapestry of forces holding atoms together, enabling us to unravel the secrets of bonding and molecular structures.\\u000b\\u000bFrom the vibrant hues of blooming flowers to the tantalizing flavors of our favorite foods, chemistry orchestrates the symphony of life. It underpins the intricate mechanisms of metabolism, the process by which living organisms convert energy from food into usable forms. Understanding chemical reactions illuminates the intricate interplay between our bodies and the environment, revealing the profound impact of nutrition, pharmaceuticals, and various environmental exposures on our health and well-being.\\u000b\\u000bThe discoveries of chemistry have revolutionized the way we live. From the transformative power of electricity to the development of innovative materials like plastics and semiconductors, chemistry has ignited a wave of technological advancements that have shaped modern society. By manipulating the properties of matter, chemists have synthesized an array of materials with tailored properties, paving the way for breakthroughs in industries ranging from electronics to medicine. The field continues to push the boundaries of knowledge, unlocking new frontiers in energy storage, sustainable technologies, and space exploration.\", \"Replace\");\n  await context.sync();\n}\n\n// 6. Summary paragraph (index 6) - full replacement, preserves run formatting\n{\n  const paras = context.document.body.paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n  const p = paras.items[6];\n  const full = p.getRange();\n  full.insertText(\"Chemistry unveils the intricacies of matter at the atomic and molecular levels, providing a foundation for understanding the behavior of substances and the forces that govern their transformations. It elucidates the processes underlying life, revealing the mechanisms of metabolism and the impact of nutrition, pharmaceuticals, and environmental factors on our health. Chemistry has played a pivotal role in technological advancements, leading to the development of innovative materials and transformative technologies that have shaped modern society. Its ongoing discoveries hold the promise of addressing global challenges and shaping a sustainable future. Chemistry stands as a testament to the power of science to unravel the mysteries of the universe and improve the human condition.\", \"Replace\");\n  await context.sync();\n}\n\n// 7. Append a new empty paragraph at the end of the body\n{\n  const body = context.document.body;\n  body.insertParagraph(\"\", \"End\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Line-break character used inside paragraphs that contain manual breaks ---\n$BR = [char]11\n\n# 1. Title\n$find = $d.Content.Find\n$find.Text = 'Galaxies: Cosmic Tapestry of Wonders'\n$find.Replacement.Text = 'Exploring Chemistry - The Building Blocks of Life'\n$find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2)\n\n# 2. Author\n$find = $d.Content.Find\n$find.Text = 'Katie Leonard'\n$find.Replacement.Text = 'Maya Singh'\n$find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2)\n\n# 3. Email local-part\n$find = $d.Content.Find\n$find.Text = 'intelligentwriter349@digitalworld'\n$find.Replacement.Text = 'mrssingh@sunnyvalehigh'\n$find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2)\n\n# 4. Email domain suffix\n$find = $d.Content.Find\n$find.Text = 'tech'\n$find.Replacement.Text = 'edu'\n$find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2)\n\n# 5. Main body paragraph (paragraph 5, with manual line breaks) - full text replacement\n$p4 = 'Chemistry is often regarded as the central science due to its profound influence on various fields, spanning biology, medicine, and materials science.'\n$p4 += ' Its study unveils the fundamental principles governing the behavior of matter at the atomic and molecular levels.'\n$p4 += ' By delving into the intricate dance of chemical reactions, we uncover the secrets behind the formation and transformation of substances, unlocking the enigmas that shape our world.'\n$p4 += ' Chemistry reveals the tapestry of forces holding atoms together, enabling us to unravel the secrets of bonding and molecular structures.'\n$p4 += $BR\n$p4 += $BR\n$p4 += 'From the vibrant hues of blooming flowers to the tantalizing flavors of our favorite foods, chemistry orchestrates the symphony of life.'\n$p4 += ' It underpins the intricate mechanisms of metabolism, the process by which living organisms convert energy from food into usable forms.'\n$p4 += ' Understanding chemical reactions illuminates the intricate interplay between our bodies and the environment, revealing the profound impact of nutrition, pharmaceuticals, and various environmental exposures on our health and well-being.'\n$p4 += $BR\n$p4 += $BR\n$p4 += 'The discoveries of chemistry have revolutionized the way we live.'\n$p4 += ' From the transformative power of electricity to the development of innovative materials like plastics and semiconductors, chemistry has ignited a wave of technological advancements that have shaped modern society.'\n$p4 += ' By manipulating the properties of matter, chemists have synthesized an array of materials with tailored properties, paving the way for breakthroughs in industries ranging from electronics to medicine.'\n$p4 += ' The field continues to push the boundaries of knowledge, unlocking new frontiers in energy storage, sustainable technologies, and space exploration.'\n\n$mainPara = $d.Paragraphs(5)\n$r = $mainPara.Range\n$r.End = $r.End - 1\n$r.Text = $p4\n\n# 6. Summary paragraph (last paragraph) - full text replacement\n$summary = 'Chemistry unveils the intricacies of matter at the atomic and molecular levels, providing a foundation for understanding the behavior of substances and the forces that govern their transformations.'\n$summary += ' It elucidates the processes underlying life, revealing the mechanisms of metabolism and the impact of nutrition, pharmaceuticals, and environmental factors on our health.'\n$summary += ' Chemistry has played a pivotal role in technological advancements, leading to the development of innovative materials and transformative technologies that have shaped modern society.'\n$summary += ' Its ongoing discoveries hold the promise of addressing global challenges and shaping a sustainable future.'\n$summary += ' Chemistry stands as a testament to the power of science to unravel the mysteries of the universe and improve the human condition.'\n\n$count = $d.Paragraphs.Count\n$summaryPara = $d.Paragraphs($count)\n$r2 = $summaryPara.Range\n$r2.End = $r2.End - 1\n$r2.Text = $summary\n\n# 7. Append a new empty paragraph at the very end of the document\n$endRange = $d.Content\n$endRange.Collapse(0)\n$endRange.InsertParagraphAfter()\n"}
